$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.69214141368866
$ws.Range("B1").Value = 2.897798299789429
$ws.Range("C1").Value = 3.60251522064209
$ws.Range("D1").Value = 1.397141814231873
$ws.Range("E1").Value = 0.9340725541114807
